# Update "想去人数" (number of people interested) counts in the
# 展览 (sheet1) and 全部类型 (sheet4) worksheets to match newly scraped
# data (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 15356
$ws1.Range("F6").Value  = 416
$ws1.Range("F8").Value  = 691
$ws1.Range("F9").Value  = 15336
$ws1.Range("F11").Value = 8919
$ws1.Range("F23").Value = 58
$ws1.Range("F24").Value = 1103
$ws1.Range("F31").Value = 49
$ws1.Range("F33").Value = 238
$ws1.Range("F34").Value = 296
$ws1.Range("F37").Value = 5456

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 15356
$ws4.Range("F6").Value  = 416
$ws4.Range("F8").Value  = 691
$ws4.Range("F9").Value  = 15336
$ws4.Range("F11").Value = 8919
$ws4.Range("F24").Value = 58
$ws4.Range("F25").Value = 1103
$ws4.Range("F34").Value = 49
$ws4.Range("F36").Value = 238
$ws4.Range("F37").Value = 296
$ws4.Range("F40").Value = 5457
